$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 1.627606790236313, 0.778006192578965),
    @(3, 5.420210225960214, 1.698254007829186),
    @(4, 6.591740668470884, 2.575414728183344),
    @(5, 7.182610763386309, 3.420968031547837),
    @(6, 9.54894408440803, 4.291343091848466),
    @(7, 9.809062539082191, 5.14912159066506),
    @(8, 10.86023555894822, 6.064601260393756),
    @(9, 10.96363787363942, 6.910068113998847),
    @(10, 13.19122909797271, 7.795044213838056),
    @(11, 13.28426613040677, 8.714102908051263),
    @(12, 15.63808196947514, 9.653178775928811),
    @(13, 16.45199117475392, 10.47085406642845),
    @(14, 16.59506908415309, 11.5321921106392),
    @(15, 24.04242650802308, 12.2707338107633),
    @(16, 24.53452056476805, 13.05911550679832),
    @(17, 24.63242927057482, 14.04347764877449),
    @(18, 36.75854155028383, 15.09413320731399),
    @(19, 39.19492893040838, 15.9175725765831),
    @(20, 39.2373260793383, 16.60216156144018),
    @(21, 40.020837059205, 17.38110114287156),
    @(22, 42.32416030474977, 18.30352728202346),
    @(23, 42.40376079154596, 19.2007761508452),
    @(24, 43.48793803060114, 20.01574244671098),
    @(25, 43.66174884861912, 21.01192659001715),
    @(26, 46.81527741613953, 21.88634908362241),
    @(27, 48.53768042923998, 22.66759313974137),
    @(28, 48.74450444027331, 23.79051347396079),
    @(29, 49.30353739840233, 24.52521177786977),
    @(30, 49.35587207012214, 25.42096518139907),
    @(31, 53.9226453190363, 26.34027674056534),
    @(32, 54.95462682317561, 27.15266366180861),
    @(33, 54.97349058826109, 27.98932613038466),
    @(34, 55.00129252319459, 28.96672671057427),
    @(35, 57.16163754720429, 29.93279202655354),
    @(36, 61.07680374657252, 30.76771385763023),
    @(37, 64.06403161474454, 31.54182208199699),
    @(38, 67.64356679356021, 32.28413884902747),
    @(39, 67.70376005289809, 33.07160212939309),
    @(40, 68.71768602056122, 34.03230933519956),
    @(41, 68.80326782417994, 34.98637856527873),
    @(42, 77.93435391790203, 35.88753197476004),
    @(43, 78.05869032460444, 36.77767494557679),
    @(44, 82.21227588852798, 37.62463906923904),
    @(45, 89.07080627576225, 38.48206823317634),
    @(46, 95.33051411874965, 39.26373160723004),
    @(47, 95.40295807927988, 40.13965404203887),
    @(48, 96.39327355045396, 40.86601089634311),
    @(49, 96.45968427560688, 41.80237025613934),
    @(50, 96.52122777031076, 42.74119348804561)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}